# Atualização de bases das ligas, do dia: 23-02-2024 às 23:34
# The odds-feed re-sync reshuffled which fixture occupies which row for a
# handful of matches. We re-read the current (pre-edit) values for every
# affected row and then re-distribute them according to the new row
# ordering, leaving column A (the sequential row id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC hold all of the per-match data (id, teams, odds, etc.);
# column A is just the row's sequence number and must stay put.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowData($rowNum) {
    $row = @{}
    foreach ($col in $cols) {
        $addr = "$col$rowNum"
        $row[$col] = $ws.Range($addr).Value()
    }
    return $row
}

function Set-RowData($rowNum, $row) {
    foreach ($col in $cols) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $row[$col]
    }
}

# Snapshot the current contents of every affected row before overwriting
# anything (several rows feed each other, including a 3-way rotation).
$row227 = Get-RowData 227
$row228 = Get-RowData 228
$row230 = Get-RowData 230
$row231 = Get-RowData 231
$row233 = Get-RowData 233
$row238 = Get-RowData 238
$row241 = Get-RowData 241

# Row 227 <-> Row 228 (simple swap)
Set-RowData 227 $row228
Set-RowData 228 $row227

# Row 230 -> Row 231 -> Row 233 -> Row 230 (3-way rotation)
Set-RowData 230 $row231
Set-RowData 231 $row233
Set-RowData 233 $row230

# Row 238 <-> Row 241 (simple swap)
Set-RowData 238 $row241
Set-RowData 241 $row238
